# Generate Report for Handback
# Populate the "Latest Target File" / "Latest Handback File" / "Latest Handback DateTime" /
# "Error Detail" columns for the bbdea7cb-8918-44b7-918e-fb00e6c81f68 row (row 6) on both the
# zh-cn and de-de localization-status sheets.

$wb = $excel.ActiveWorkbook

$targetUrl = "https://github.com/OpenLocalizationTestOrg/ol-test4/blob/494e729a0dc648da605af2759ef403856318cda4/e2e/bbdea7cb-8918-44b7-918e-fb00e6c81f68.md"
$targetDisplay = "bbdea7cb-8918-44b7-918e-fb00e6c81f68.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test4/blob/a2994b5124b31d56cbc9145f18983b9ceea72dad/e2e/bbdea7cb-8918-44b7-918e-fb00e6c81f68.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test4/blob/494e729a0dc648da605af2759ef403856318cda4/e2e/bbdea7cb-8918-44b7-918e-fb00e6c81f68.md."

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")

# J6: Latest Target File -> hyperlink to the handback markdown file
$wsZh.Hyperlinks.Add($wsZh.Cells.Item(6, 10), $targetUrl, [Type]::Missing, [Type]::Missing, $targetDisplay)

# K6: Latest Handback File
$wsZh.Cells.Item(6, 11).Value = "bbdea7cb-8918-44b7-918e-fb00e6c81f68.6ad35c227a43d1831b05370aa5a076793e028f1d.zh-cn.xlf"

# L6: Latest Handback DateTime
$wsZh.Cells.Item(6, 12).Value = "2017-02-21 03:44:51"

# R6: Error Detail
$wsZh.Cells.Item(6, 18).Value = $errorDetail

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")

# J6: Latest Target File -> hyperlink to the handback markdown file
$wsDe.Hyperlinks.Add($wsDe.Cells.Item(6, 10), $targetUrl, [Type]::Missing, [Type]::Missing, $targetDisplay)

# K6: Latest Handback File
$wsDe.Cells.Item(6, 11).Value = "bbdea7cb-8918-44b7-918e-fb00e6c81f68.6ad35c227a43d1831b05370aa5a076793e028f1d.de-de.xlf"

# L6: Latest Handback DateTime
$wsDe.Cells.Item(6, 12).Value = "2017-02-21 03:45:14"

# R6: Error Detail
$wsDe.Cells.Item(6, 18).Value = $errorDetail
